$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (Q) that mirrors the formatting of the existing
# "2019" column (P), then overwrite the header value for the new year.
$ws.Range("P4").Copy($ws.Range("Q4")) | Out-Null
$ws.Range("Q4").Value = 2020

# The data rows for 2020 repeat the same figures as 2019, so copy the whole
# cell (value + style) straight across.
$ws.Range("P5").Copy($ws.Range("Q5")) | Out-Null
$ws.Range("P6").Copy($ws.Range("Q6")) | Out-Null
$ws.Range("P7").Copy($ws.Range("Q7")) | Out-Null
$ws.Range("P8").Copy($ws.Range("Q8")) | Out-Null

# Restore the active cell/selection that was active when the file was saved.
$ws.Range("O12").Select() | Out-Null
